# "programação das figuras do grupo 16"
# - Column C (Ano) dates change from 31/12/YYYY to 01/01/YYYY, keeping them
#   as literal text (not Excel date serials).
# - Header row (A1:D1) gains a thin border all around plus vertical-top
#   alignment (in addition to the existing horizontal-center alignment).
# - Page margins reset to Excel's "normal" defaults.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column C to text format BEFORE writing the new values so the
# "01/01/2016"-style strings are not auto-converted into date serials
# (Excel would otherwise interpret a valid dd/mm/yyyy string as a date).
$dates = $ws.Range("C2:C19")
$dates.NumberFormat = "@"

$ws.Range("C2").Value = "01/01/2016"
$ws.Range("C3").Value = "01/01/2017"
$ws.Range("C4").Value = "01/01/2018"
$ws.Range("C5").Value = "01/01/2019"
$ws.Range("C6").Value = "01/01/2022"
$ws.Range("C7").Value = "01/01/2023"
$ws.Range("C8").Value = "01/01/2016"
$ws.Range("C9").Value = "01/01/2017"
$ws.Range("C10").Value = "01/01/2018"
$ws.Range("C11").Value = "01/01/2019"
$ws.Range("C12").Value = "01/01/2022"
$ws.Range("C13").Value = "01/01/2023"
$ws.Range("C14").Value = "01/01/2016"
$ws.Range("C15").Value = "01/01/2017"
$ws.Range("C16").Value = "01/01/2018"
$ws.Range("C17").Value = "01/01/2019"
$ws.Range("C18").Value = "01/01/2022"
$ws.Range("C19").Value = "01/01/2023"

# Header row: add a thin border around each cell and top-align the (already
# horizontally centered) bold header text.
$header = $ws.Range("A1:D1")
$header.Borders.LineStyle = 1
$header.VerticalAlignment = -4160

# Restore Excel's default page margins (0.75"/0.75"/1"/1"/0.5"/0.5").
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36
